# descw-1252 finalize tab38 rpt
# Update the "fiscal year" template placeholder used in the report model
# from `{#fy=d.fiscal_year}` to `{#fy=d.fiscal}` (B11 on Sheet1), and
# leave the workbook's selection where the author left it (B16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "{#fy=d.fiscal}"

$ws.Range("B16").Select()
